$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.004.63'
$ws.Range("E2").Value = '  +2.41%  '

$ws.Range("D3").Value = '3.386.03'
$ws.Range("E3").Value = '  +2.11%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '559.20'
$ws.Range("E5").Value = '  +2.37%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.10'
$ws.Range("E6").Value = '  +1.21%  '

$ws.Range("E7").Value = '  +1.51%  '

$ws.Range("D8").Value = '3.377.23'
$ws.Range("E8").Value = '  +2.12%  '

$ws.Range("E9").Value = '  +0.00%  '

$ws.Range("E10").Value = '  +11.98%  '

$ws.Range("E11").Value = '  +3.57%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.23'
$ws.Range("E12").Value = '  +2.43%  '

$ws.Range("E13").Value = '  +5.85%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.09'
$ws.Range("E14").Value = '  +3.01%  '

$ws.Range("D15").Value = '3.927.58'
$ws.Range("E15").Value = '  +2.10%  '

$ws.Range("E16").Value = '  +1.51%  '

$ws.Range("E17").Value = '  +2.38%  '

$ws.Range("D18").Value = '3.384.60'
$ws.Range("E18").Value = '  +2.15%  '

$ws.Range("D19").Value = '64.960.39'
$ws.Range("E19").Value = '  +2.33%  '

$ws.Range("E20").Value = '  +1.89%  '

$ws.Range("E21").Value = '  +2.28%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '470.72'
$ws.Range("E22").Value = '  +14.33%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.95'
$ws.Range("E23").Value = '  +13.15%  '

$ws.Range("E24").Value = '  +2.76%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.85'
$ws.Range("E25").Value = '  +5.06%  '

$ws.Range("E26").Value = '  -0.69%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.89'
$ws.Range("E27").Value = '  +6.81%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.76'
$ws.Range("E28").Value = '  +2.67%  '

$ws.Range("E29").Value = '  +2.38%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.83'
$ws.Range("E30").Value = '  +6.55%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.67'
$ws.Range("E31").Value = '  +5.35%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.49'
$ws.Range("E32").Value = '  +1.78%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '572.17'
$ws.Range("E33").Value = '  -0.69%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '61.35'
$ws.Range("E34").Value = '  +6.51%  '

$ws.Range("E35").Value = '  +2.28%  '

$ws.Range("E36").Value = '  +0.00%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.61'
$ws.Range("E37").Value = '  +6.58%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.140'
$ws.Range("E38").Value = '  -4.67%  '

$ws.Range("E39").Value = '  +2.29%  '

$ws.Range("D40").Value = '0.0₃0750'
$ws.Range("E40").Value = '  +2.47%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.369'
$ws.Range("E41").Value = '  +1.51%  '

$ws.Range("D42").Value = '3.088.78'
$ws.Range("E42").Value = '  -0.75%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  +0.00%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.84'
$ws.Range("E44").Value = '  +3.18%  '

$ws.Range("E45").Value = '  +4.31%  '

$ws.Range("E46").Value = '  +5.65%  '

$ws.Range("E47").Value = '  +2.72%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.13'
$ws.Range("E48").Value = '  -2.59%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.59'
$ws.Range("E49").Value = '  +0.02%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '139.38'
$ws.Range("E50").Value = '  +5.35%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.26'
$ws.Range("E51").Value = '  +3.77%  '
